# Scheduled runner: refresh Market Board price snapshots in the Chocobo
# profit-tracking sheets (currentAveragePrice / LevePrice / LeveProfit
# columns H:N) for a handful of leves across ALC, ARM, BSM, CRP, GSM, WVR.

$wb = $excel.ActiveWorkbook

function Set-Row {
    param(
        [string]$SheetName,
        [int]$Row,
        [hashtable]$Values   # column letter -> new value (or $null to clear)
    )
    $ws = $wb.Worksheets.Item($SheetName)
    foreach ($col in $Values.Keys) {
        $cell = $ws.Range("$col$Row")
        $val = $Values[$col]
        if ($null -eq $val) {
            $cell.ClearContents()
        } else {
            $cell.Value2 = $val
        }
    }
}

# ---------------- ALC ----------------
Set-Row "ALC" 123 @{ H = 0; J = 0; L = 0; N = $null }
Set-Row "ALC" 137 @{ H = 2298.2444; I = 1539.919; J = 5805.5; K = 4619.757000000001; L = 17416.5; M = -2069.757000000001; N = -22516.5 }
Set-Row "ALC" 138 @{ H = 2251.68; I = 1043.1333; J = 2769.6287; K = 3129.3999; L = 8308.8861; M = 2010.6001; N = -18588.8861 }

# ---------------- ARM ----------------
Set-Row "ARM" 26  @{ H = 5000; J = 0; L = 0; N = $null }
Set-Row "ARM" 32  @{ H = 9960.987999999999; I = 6542.4653; K = 6542.4653; M = -6255.4653 }
Set-Row "ARM" 61  @{ H = 1537.7858; I = 968.9; J = 2960; K = 968.9; L = 2960; M = -756.9; N = -3384 }
Set-Row "ARM" 64  @{ H = 48993.57; J = 48993.57; L = 48993.57; N = -49489.57 }
Set-Row "ARM" 67  @{ H = 48993.57; J = 48993.57; L = 48993.57; N = -50709.57 }
Set-Row "ARM" 74  @{ H = 1342.4138; I = 989.82355; K = 989.82355; M = -115.82355 }
Set-Row "ARM" 77  @{ H = 1342.4138; I = 989.82355; K = 4949.117749999999; M = -581.1177499999994 }
Set-Row "ARM" 135 @{ H = 47003; J = 47003; L = 47003; N = -57143 }
Set-Row "ARM" 136 @{ H = 1537.7858; I = 968.9; J = 2960; K = 2906.7; L = 8880; M = -356.6999999999998; N = -13980 }

# ---------------- BSM ----------------
Set-Row "BSM" 62  @{ H = 50181; J = 50181; L = 50181; N = -51553 }
Set-Row "BSM" 65  @{ H = 50181; J = 50181; L = 150543; N = -157407 }
Set-Row "BSM" 137 @{ H = 34817.6; J = 34817.6; L = 34817.6; N = -45017.6 }

# ---------------- CRP ----------------
Set-Row "CRP" 31  @{ H = 2783.0425; I = 1271.8485; J = 6345.143; K = 1271.8485; L = 6345.143; M = -976.8485000000001; N = -6935.143 }
Set-Row "CRP" 32  @{ H = 20000; I = 20000; K = 20000; M = -19684 }
Set-Row "CRP" 34  @{ H = 2783.0425; I = 1271.8485; J = 6345.143; K = 1271.8485; L = 6345.143; M = -1069.8485; N = -6749.143 }
Set-Row "CRP" 58  @{ H = 1981.5374; I = 1669.8853; J = 5150; K = 1669.8853; L = 5150; M = -1466.8853; N = -5556 }
Set-Row "CRP" 68  @{ H = 47676.1; J = 47676.1; L = 47676.1; N = -49174.1 }
Set-Row "CRP" 71  @{ H = 47676.1; J = 47676.1; L = 143028.3; N = -150516.3 }
Set-Row "CRP" 125 @{ H = 35215; J = 35215; L = 35215; N = -40135 }
Set-Row "CRP" 132 @{ H = 3257.0881; I = 3016.95; J = 3600.1428; K = 9050.849999999999; L = 10800.4284; M = -6520.849999999999; N = -15860.4284 }
Set-Row "CRP" 136 @{ H = 1981.5374; I = 1669.8853; J = 5150; K = 5009.6559; L = 15450; M = -2459.6559; N = -20550 }

# ---------------- GSM ----------------
Set-Row "GSM" 31  @{ H = 20765.5; I = 11531; J = 30000; K = 11531; L = 30000; M = -11239; N = -30584 }
Set-Row "GSM" 37  @{ H = 20765.5; I = 11531; J = 30000; K = 11531; L = 30000; M = -11254; N = -30554 }
Set-Row "GSM" 70  @{ H = 6381.0356; I = 5741.381; K = 5741.381; M = -5471.381 }
Set-Row "GSM" 73  @{ H = 6381.0356; I = 5741.381; K = 5741.381; M = -4805.381 }

# ---------------- WVR ----------------
Set-Row "WVR" 40  @{ H = 29499; I = 0; J = 29499; K = 0; L = 29499; M = $null; N = -29797 }
Set-Row "WVR" 132 @{ H = 6948372.5; I = 5367.273; J = 12823223; K = 16101.819; L = 38469669; M = -13571.819; N = -38474729 }
